$d = $word.ActiveDocument

$d.Content.Find.Execute("20+0=", $true, $false, $false, $false, $false, $true, 1, $false, "65-51=", 2) | Out-Null
$d.Content.Find.Execute("15+41=", $true, $false, $false, $false, $false, $true, 1, $false, "20+20=", 2) | Out-Null
$d.Content.Find.Execute("78+14=", $true, $false, $false, $false, $false, $true, 1, $false, "86-73=", 2) | Out-Null
$d.Content.Find.Execute("31+26=", $true, $false, $false, $false, $false, $true, 1, $false, "69-68=", 2) | Out-Null
$d.Content.Find.Execute("18+76=", $true, $false, $false, $false, $false, $true, 1, $false, "85-19=", 2) | Out-Null
$d.Content.Find.Execute("80-27=", $true, $false, $false, $false, $false, $true, 1, $false, "54+41=", 2) | Out-Null
$d.Content.Find.Execute("99-90=", $true, $false, $false, $false, $false, $true, 1, $false, "29-28=", 2) | Out-Null
$d.Content.Find.Execute("68+21=", $true, $false, $false, $false, $false, $true, 1, $false, "13-6=", 2) | Out-Null
$d.Content.Find.Execute("87-72=", $true, $false, $false, $false, $false, $true, 1, $false, "44-38=", 2) | Out-Null
$d.Content.Find.Execute("84-35=", $true, $false, $false, $false, $false, $true, 1, $false, "69-42=", 2) | Out-Null
$d.Content.Find.Execute("69-15=", $true, $false, $false, $false, $false, $true, 1, $false, "50-18=", 2) | Out-Null
$d.Content.Find.Execute("90-84=", $true, $false, $false, $false, $false, $true, 1, $false, "79-61=", 2) | Out-Null
$d.Content.Find.Execute("8+53=", $true, $false, $false, $false, $false, $true, 1, $false, "31-27=", 2) | Out-Null
$d.Content.Find.Execute("48+37=", $true, $false, $false, $false, $false, $true, 1, $false, "71-49=", 2) | Out-Null
$d.Content.Find.Execute("10+33=", $true, $false, $false, $false, $false, $true, 1, $false, "39-13=", 2) | Out-Null
$d.Content.Find.Execute("73+23=", $true, $false, $false, $false, $false, $true, 1, $false, "47+13=", 2) | Out-Null
$d.Content.Find.Execute("40+32=", $true, $false, $false, $false, $false, $true, 1, $false, "94-4=", 2) | Out-Null
$d.Content.Find.Execute("59-18=", $true, $false, $false, $false, $false, $true, 1, $false, "9+84=", 2) | Out-Null
$d.Content.Find.Execute("12+23=", $true, $false, $false, $false, $false, $true, 1, $false, "58-19=", 2) | Out-Null
$d.Content.Find.Execute("79-15=", $true, $false, $false, $false, $false, $true, 1, $false, "36-25=", 2) | Out-Null
$d.Content.Find.Execute("93-92=", $true, $false, $false, $false, $false, $true, 1, $false, "56+11=", 2) | Out-Null
$d.Content.Find.Execute("29+1=", $true, $false, $false, $false, $false, $true, 1, $false, "88-85=", 2) | Out-Null
$d.Content.Find.Execute("33-12=", $true, $false, $false, $false, $false, $true, 1, $false, "30+46=", 2) | Out-Null
$d.Content.Find.Execute("55-19=", $true, $false, $false, $false, $false, $true, 1, $false, "10+64=", 2) | Out-Null
$d.Content.Find.Execute("17+15=", $true, $false, $false, $false, $false, $true, 1, $false, "85-14=", 2) | Out-Null
$d.Content.Find.Execute("53+46=", $true, $false, $false, $false, $false, $true, 1, $false, "36+28=", 2) | Out-Null
$d.Content.Find.Execute("79-76=", $true, $false, $false, $false, $false, $true, 1, $false, "33-5=", 2) | Out-Null
$d.Content.Find.Execute("83+8=", $true, $false, $false, $false, $false, $true, 1, $false, "68-22=", 2) | Out-Null
$d.Content.Find.Execute("74-69=", $true, $false, $false, $false, $false, $true, 1, $false, "96-46=", 2) | Out-Null
$d.Content.Find.Execute("21+71=", $true, $false, $false, $false, $false, $true, 1, $false, "86-51=", 2) | Out-Null
$d.Content.Find.Execute("43+27=", $true, $false, $false, $false, $false, $true, 1, $false, "70-47=", 2) | Out-Null
$d.Content.Find.Execute("67-6=", $true, $false, $false, $false, $false, $true, 1, $false, "6+11=", 2) | Out-Null
$d.Content.Find.Execute("79-8=", $true, $false, $false, $false, $false, $true, 1, $false, "19-3=", 2) | Out-Null
$d.Content.Find.Execute("18-17=", $true, $false, $false, $false, $false, $true, 1, $false, "12+52=", 2) | Out-Null
$d.Content.Find.Execute("29-14=", $true, $false, $false, $false, $false, $true, 1, $false, "21+33=", 2) | Out-Null
$d.Content.Find.Execute("28+26=", $true, $false, $false, $false, $false, $true, 1, $false, "59-25=", 2) | Out-Null
$d.Content.Find.Execute("46+11=", $true, $false, $false, $false, $false, $true, 1, $false, "8+29=", 2) | Out-Null
$d.Content.Find.Execute("45-19=", $true, $false, $false, $false, $false, $true, 1, $false, "18+65=", 2) | Out-Null
$d.Content.Find.Execute("49-46=", $true, $false, $false, $false, $false, $true, 1, $false, "35+33=", 2) | Out-Null
$d.Content.Find.Execute("0+5=", $true, $false, $false, $false, $false, $true, 1, $false, "57+38=", 2) | Out-Null
$d.Content.Find.Execute("75-74=", $true, $false, $false, $false, $false, $true, 1, $false, "73-33=", 2) | Out-Null
$d.Content.Find.Execute("76-29=", $true, $false, $false, $false, $false, $true, 1, $false, "32-14=", 2) | Out-Null
$d.Content.Find.Execute("0+79=", $true, $false, $false, $false, $false, $true, 1, $false, "62-11=", 2) | Out-Null
$d.Content.Find.Execute("82-12=", $true, $false, $false, $false, $false, $true, 1, $false, "49-29=", 2) | Out-Null
$d.Content.Find.Execute("21+10=", $true, $false, $false, $false, $false, $true, 1, $false, "16-3=", 2) | Out-Null
$d.Content.Find.Execute("26+1=", $true, $false, $false, $false, $false, $true, 1, $false, "93-4=", 2) | Out-Null
$d.Content.Find.Execute("49+0=", $true, $false, $false, $false, $false, $true, 1, $false, "53+24=", 2) | Out-Null
$d.Content.Find.Execute("20-3=", $true, $false, $false, $false, $false, $true, 1, $false, "30-20=", 2) | Out-Null
$d.Content.Find.Execute("2+79=", $true, $false, $false, $false, $false, $true, 1, $false, "64-2=", 2) | Out-Null
$d.Content.Find.Execute("95-12=", $true, $false, $false, $false, $false, $true, 1, $false, "3+56=", 2) | Out-Null
$d.Content.Find.Execute("48+51=", $true, $false, $false, $false, $false, $true, 1, $false, "34+13=", 2) | Out-Null
$d.Content.Find.Execute("20+25=", $true, $false, $false, $false, $false, $true, 1, $false, "51+8=", 2) | Out-Null
$d.Content.Find.Execute("30+42=", $true, $false, $false, $false, $false, $true, 1, $false, "95-28=", 2) | Out-Null
$d.Content.Find.Execute("17-9=", $true, $false, $false, $false, $false, $true, 1, $false, "56-51=", 2) | Out-Null
$d.Content.Find.Execute("92+0=", $true, $false, $false, $false, $false, $true, 1, $false, "96-26=", 2) | Out-Null
$d.Content.Find.Execute("62+4=", $true, $false, $false, $false, $false, $true, 1, $false, "97-31=", 2) | Out-Null
$d.Content.Find.Execute("27+63=", $true, $false, $false, $false, $false, $true, 1, $false, "58+31=", 2) | Out-Null
$d.Content.Find.Execute("16+78=", $true, $false, $false, $false, $false, $true, 1, $false, "10+16=", 2) | Out-Null
$d.Content.Find.Execute("45+2=", $true, $false, $false, $false, $false, $true, 1, $false, "8+74=", 2) | Out-Null
$d.Content.Find.Execute("42+33=", $true, $false, $false, $false, $false, $true, 1, $false, "92-8=", 2) | Out-Null
$d.Content.Find.Execute("2+75=", $true, $false, $false, $false, $false, $true, 1, $false, "14-1=", 2) | Out-Null
$d.Content.Find.Execute("84-21=", $true, $false, $false, $false, $false, $true, 1, $false, "17+35=", 2) | Out-Null
$d.Content.Find.Execute("67-27=", $true, $false, $false, $false, $false, $true, 1, $false, "89-4=", 2) | Out-Null
$d.Content.Find.Execute("91-69=", $true, $false, $false, $false, $false, $true, 1, $false, "29-23=", 2) | Out-Null
$d.Content.Find.Execute("78+10=", $true, $false, $false, $false, $false, $true, 1, $false, "88-57=", 2) | Out-Null
$d.Content.Find.Execute("46+12=", $true, $false, $false, $false, $false, $true, 1, $false, "43+0=", 2) | Out-Null
$d.Content.Find.Execute("38-24=", $true, $false, $false, $false, $false, $true, 1, $false, "14-9=", 2) | Out-Null
$d.Content.Find.Execute("73-59=", $true, $false, $false, $false, $false, $true, 1, $false, "29+31=", 2) | Out-Null
$d.Content.Find.Execute("32+21=", $true, $false, $false, $false, $false, $true, 1, $false, "53-13=", 2) | Out-Null
$d.Content.Find.Execute("48+43=", $true, $false, $false, $false, $false, $true, 1, $false, "92-47=", 2) | Out-Null
$d.Content.Find.Execute("14+60=", $true, $false, $false, $false, $false, $true, 1, $false, "87-55=", 2) | Out-Null
$d.Content.Find.Execute("79+10=", $true, $false, $false, $false, $false, $true, 1, $false, "70-39=", 2) | Out-Null
$d.Content.Find.Execute("77-70=", $true, $false, $false, $false, $false, $true, 1, $false, "48-45=", 2) | Out-Null
$d.Content.Find.Execute("97-90=", $true, $false, $false, $false, $false, $true, 1, $false, "78-69=", 2) | Out-Null
$d.Content.Find.Execute("65-49=", $true, $false, $false, $false, $false, $true, 1, $false, "1+73=", 2) | Out-Null
$d.Content.Find.Execute("26+37=", $true, $false, $false, $false, $false, $true, 1, $false, "83+11=", 2) | Out-Null
$d.Content.Find.Execute("77-46=", $true, $false, $false, $false, $false, $true, 1, $false, "39+21=", 2) | Out-Null
$d.Content.Find.Execute("66-36=", $true, $false, $false, $false, $false, $true, 1, $false, "69-0=", 2) | Out-Null
$d.Content.Find.Execute("20+42=", $true, $false, $false, $false, $false, $true, 1, $false, "10+57=", 2) | Out-Null
$d.Content.Find.Execute("86-44=", $true, $false, $false, $false, $false, $true, 1, $false, "14+40=", 2) | Out-Null
$d.Content.Find.Execute("67-28=", $true, $false, $false, $false, $false, $true, 1, $false, "2+2=", 2) | Out-Null
$d.Content.Find.Execute("35+46=", $true, $false, $false, $false, $false, $true, 1, $false, "59-42=", 2) | Out-Null
$d.Content.Find.Execute("55-26=", $true, $false, $false, $false, $false, $true, 1, $false, "20+22=", 2) | Out-Null
$d.Content.Find.Execute("85-81=", $true, $false, $false, $false, $false, $true, 1, $false, "86-81=", 2) | Out-Null
$d.Content.Find.Execute("23+5=", $true, $false, $false, $false, $false, $true, 1, $false, "43-38=", 2) | Out-Null
$d.Content.Find.Execute("93-5=", $true, $false, $false, $false, $false, $true, 1, $false, "21-7=", 2) | Out-Null
$d.Content.Find.Execute("55+2=", $true, $false, $false, $false, $false, $true, 1, $false, "81+10=", 2) | Out-Null
$d.Content.Find.Execute("14+39=", $true, $false, $false, $false, $false, $true, 1, $false, "98-44=", 2) | Out-Null
$d.Content.Find.Execute("37+19=", $true, $false, $false, $false, $false, $true, 1, $false, "33-28=", 2) | Out-Null
$d.Content.Find.Execute("47+40=", $true, $false, $false, $false, $false, $true, 1, $false, "16+11=", 2) | Out-Null
$d.Content.Find.Execute("91+5=", $true, $false, $false, $false, $false, $true, 1, $false, "35-12=", 2) | Out-Null
$d.Content.Find.Execute("90-59=", $true, $false, $false, $false, $false, $true, 1, $false, "52-25=", 2) | Out-Null
$d.Content.Find.Execute("33+2=", $true, $false, $false, $false, $false, $true, 1, $false, "57-50=", 2) | Out-Null
$d.Content.Find.Execute("94-59=", $true, $false, $false, $false, $false, $true, 1, $false, "85+12=", 2) | Out-Null
$d.Content.Find.Execute("81-2=", $true, $false, $false, $false, $false, $true, 1, $false, "66-53=", 2) | Out-Null
$d.Content.Find.Execute("59-28=", $true, $false, $false, $false, $false, $true, 1, $false, "90-41=", 2) | Out-Null
$d.Content.Find.Execute("20+40=", $true, $false, $false, $false, $false, $true, 1, $false, "95-54=", 2) | Out-Null
$d.Content.Find.Execute("25+20=", $true, $false, $false, $false, $false, $true, 1, $false, "64-30=", 2) | Out-Null
$d.Content.Find.Execute("90+7=", $true, $false, $false, $false, $false, $true, 1, $false, "41+37=", 2) | Out-Null
$d.Content.Find.Execute("32+66=", $true, $false, $false, $false, $false, $true, 1, $false, "76-27=", 2) | Out-Null
